# Fixed calcmem for RobinHoodHashing: correct the raw per-trial benchmark
# measurements on Sheet1 (Static Array vs RobinHood Hashing memory
# consumption samples, rows 4-9 and 14-19). Columns M/N hold AVERAGE()
# formulas and column O the M/N ratio; all three recalc automatically from
# the corrected inputs below, and chart1/chart2 plot those M/N columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("C4").Value = 2196480
$ws.Range("D4").Value = 499170
$ws.Range("E4").Value = 2197140
$ws.Range("F4").Value = 499320
$ws.Range("G4").Value = 2250380
$ws.Range("H4").Value = 511420
$ws.Range("I4").Value = 2309560
$ws.Range("J4").Value = 524870
$ws.Range("K4").Value = 2222880
$ws.Range("L4").Value = 505170

# Row 5
$ws.Range("C5").Value = 10804200
$ws.Range("D5").Value = 2455520
$ws.Range("E5").Value = 10980640
$ws.Range("F5").Value = 2495570
$ws.Range("G5").Value = 10806620
$ws.Range("H5").Value = 2456170
$ws.Range("I5").Value = 10602900
$ws.Range("J5").Value = 2409720
$ws.Range("K5").Value = 10527000
$ws.Range("L5").Value = 2392470

# Row 6
$ws.Range("C6").Value = 20887900
$ws.Range("D6").Value = 4747220
$ws.Range("E6").Value = 21405780
$ws.Range("F6").Value = 4864970
$ws.Range("G6").Value = 20948180
$ws.Range("H6").Value = 4760920
$ws.Range("I6").Value = 20834660
$ws.Range("J6").Value = 4735120
$ws.Range("K6").Value = 21066760
$ws.Range("L6").Value = 4787870

# Row 7
$ws.Range("C7").Value = 194067500
$ws.Range("D7").Value = 44106320
$ws.Range("E7").Value = 195299500
$ws.Range("F7").Value = 44386670
$ws.Range("G7").Value = 194030980
$ws.Range("H7").Value = 44098320
$ws.Range("I7").Value = 194644340
$ws.Range("J7").Value = 44238020
$ws.Range("K7").Value = 194726400
$ws.Range("L7").Value = 44256420

# Row 8
$ws.Range("C8").Value = 568936940
$ws.Range("D8").Value = 129305570
$ws.Range("E8").Value = 566506380
$ws.Range("F8").Value = 128752770
$ws.Range("G8").Value = 565739680
$ws.Range("H8").Value = 128578970
$ws.Range("I8").Value = 565921840
$ws.Range("J8").Value = 128620170
$ws.Range("K8").Value = 566620780
$ws.Range("L8").Value = 128779720

# Row 9
$ws.Range("C9").Value = 929984880
$ws.Range("D9").Value = 211363520
$ws.Range("E9").Value = 930845520
$ws.Range("F9").Value = 211559470
$ws.Range("G9").Value = 928468860
$ws.Range("H9").Value = 211019520
$ws.Range("I9").Value = 928998180
$ws.Range("J9").Value = 211139520
$ws.Range("K9").Value = 931036040
$ws.Range("L9").Value = 211602570

# Row 14
$ws.Range("C14").Value = 334840
$ws.Range("D14").Value = 76070
$ws.Range("E14").Value = 334840
$ws.Range("F14").Value = 76070
$ws.Range("G14").Value = 336380
$ws.Range("H14").Value = 76420
$ws.Range("I14").Value = 335720
$ws.Range("J14").Value = 76270
$ws.Range("K14").Value = 340780
$ws.Range("L14").Value = 77420

# Row 15
$ws.Range("C15").Value = 3302640
$ws.Range("D15").Value = 750570
$ws.Range("E15").Value = 3309680
$ws.Range("F15").Value = 752220
$ws.Range("G15").Value = 3308140
$ws.Range("H15").Value = 751820
$ws.Range("I15").Value = 3302640
$ws.Range("J15").Value = 750570
$ws.Range("K15").Value = 3306160
$ws.Range("L15").Value = 751370

# Row 16
$ws.Range("C16").Value = 12809280
$ws.Range("D16").Value = 2911220
$ws.Range("E16").Value = 12802900
$ws.Range("F16").Value = 2909820
$ws.Range("G16").Value = 12808180
$ws.Range("H16").Value = 2911120
$ws.Range("I16").Value = 12811480
$ws.Range("J16").Value = 2911820
$ws.Range("K16").Value = 12805100
$ws.Range("L16").Value = 2910420

# Row 17
$ws.Range("C17").Value = 155650000
$ws.Range("D17").Value = 35375570
$ws.Range("E17").Value = 155700160
$ws.Range("F17").Value = 35386620
$ws.Range("G17").Value = 155669800
$ws.Range("H17").Value = 35379770
$ws.Range("I17").Value = 155677500
$ws.Range("J17").Value = 35381520
$ws.Range("K17").Value = 155656380
$ws.Range("L17").Value = 35376970

# Row 18
$ws.Range("C18").Value = 777566240
$ws.Range("D18").Value = 176721320
$ws.Range("E18").Value = 777605180
$ws.Range("F18").Value = 176731170
$ws.Range("G18").Value = 777536980
$ws.Range("H18").Value = 176715670
$ws.Range("I18").Value = 777553480
$ws.Range("J18").Value = 176718370
$ws.Range("K18").Value = 777557880
$ws.Range("L18").Value = 176720170

# Row 19
$ws.Range("C19").Value = 1828554640
$ws.Range("D19").Value = 415586670
$ws.Range("E19").Value = 1828473460
$ws.Range("F19").Value = 415567920
$ws.Range("G19").Value = 1828524720
$ws.Range("H19").Value = 415578720
$ws.Range("I19").Value = 1828429020
$ws.Range("J19").Value = 415556570
$ws.Range("K19").Value = 1828579720
$ws.Range("L19").Value = 415592520

# Restore the active sheet/selection as last saved by the author.
[void]$ws.Activate()
[void]$ws.Range("E12:F12").Select()

